$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 ("I0") and J1 ("IF"), copying the existing header
# style (bold, bordered, centered) from H1 so the new headers match.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new I/J column data for rows 2-12.
$data = @(
    @(3, 4),
    @(8, 9),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(4, 5),
    @(6, 8),
    @(6, 6),
    @(7, 8),
    @(9, 9),
    @(7, 8)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
